# Auto-generated PowerShell Excel COM-interop script
# Applies odds updates to rows 4, 6, 8, 9, 13 per the target diff

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("G4").Value = 3.6
$ws.Range("H4").Value = 3.2
$ws.Range("I4").Value = 2.1
$ws.Range("R4").Value = 2.1
$ws.Range("S4").Value = 1.67
$ws.Range("T4").Value = 8
$ws.Range("U4").Value = 17
$ws.Range("V4").Value = 13
$ws.Range("X4").Value = 34
$ws.Range("Y4").Value = 41
$ws.Range("AB4").Value = 19
$ws.Range("AC4").Value = 67
$ws.Range("AE4").Value = 6
$ws.Range("AF4").Value = 9
$ws.Range("AH4").Value = 19
$ws.Range("AI4").Value = 21

# Row 6
$ws.Range("G6").Value = 2
$ws.Range("H6").Value = 2.92
$ws.Range("I6").Value = 4.05
$ws.Range("L6").Value = 1.45
$ws.Range("M6").Value = 2.37
$ws.Range("N6").Value = 2.27
$ws.Range("O6").Value = 1.5
$ws.Range("Q6").Value = 2.27
$ws.Range("T6").Value = 5.8
$ws.Range("U6").Value = 8.5
$ws.Range("V6").Value = 8.75
$ws.Range("W6").Value = 18
$ws.Range("X6").Value = 18.5
$ws.Range("Y6").Value = 37
$ws.Range("AA6").Value = 5.8
$ws.Range("AB6").Value = 17.5
$ws.Range("AC6").Value = 110
$ws.Range("AE6").Value = 8.5
$ws.Range("AF6").Value = 20
$ws.Range("AG6").Value = 14.5
$ws.Range("AH6").Value = 70
$ws.Range("AI6").Value = 50
$ws.Range("AJ6").Value = 65

# Row 8
$ws.Range("G8").Value = 6.5
$ws.Range("H8").Value = 3.55
$ws.Range("I8").Value = 1.55
$ws.Range("P8").Value = 1.47
$ws.Range("Q8").Value = 2.5
$ws.Range("R8").Value = 2.2
$ws.Range("S8").Value = 1.6
$ws.Range("V8").Value = 21
$ws.Range("W8").Value = 150
$ws.Range("X8").Value = 90
$ws.Range("Y8").Value = 90
$ws.Range("AA8").Value = 7.1
$ws.Range("AB8").Value = 22
$ws.Range("AF8").Value = 6.2
$ws.Range("AH8").Value = 10.5
$ws.Range("AI8").Value = 14.5

# Row 9
$ws.Range("G9").Value = 7.5
$ws.Range("H9").Value = 3.85
$ws.Range("I9").Value = 1.45
$ws.Range("K9").Value = 6.9
$ws.Range("L9").Value = 1.33
$ws.Range("M9").Value = 3.05
$ws.Range("N9").Value = 1.98
$ws.Range("O9").Value = 1.75
$ws.Range("P9").Value = 1.4
$ws.Range("Q9").Value = 2.75
$ws.Range("T9").Value = 16.5
$ws.Range("V9").Value = 23
$ws.Range("X9").Value = 90
$ws.Range("Z9").Value = 6.9
$ws.Range("AA9").Value = 7.6
$ws.Range("AE9").Value = 5.5
$ws.Range("AH9").Value = 9.5

# Row 13
$ws.Range("J13").Value = 1.02
$ws.Range("K13").Value = 21
